$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Website"
$ws.Range("A2").Value = "Tim viec"

$ws.Range("B1").Select()
